# Update the date heading.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-08-19 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-08-20 Tuesday", 2)

# Update the division problems in the table. Using Cell.Range.Text (rather
# than a global text Find/Replace) so each cell is addressed by its exact
# position — this avoids any ambiguity from duplicate/overlapping values
# that appear before and after the edits (e.g. "35÷6=5, 5").
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "51÷7=7, 2"
$t.Cell(1, 2).Range.Text  = "60÷8=7, 4"
$t.Cell(1, 3).Range.Text  = "35÷6=5, 5"
$t.Cell(1, 4).Range.Text  = "36÷4=9, 0"
$t.Cell(1, 5).Range.Text  = "81÷6=13, 3"

$t.Cell(5, 1).Range.Text  = "63÷4=15, 3"
$t.Cell(5, 2).Range.Text  = "69÷9=7, 6"
$t.Cell(5, 3).Range.Text  = "30÷7=4, 2"
$t.Cell(5, 4).Range.Text  = "45÷9=5, 0"
$t.Cell(5, 5).Range.Text  = "64÷6=10, 4"

$t.Cell(9, 1).Range.Text  = "88÷2=44, 0"
$t.Cell(9, 2).Range.Text  = "62÷8=7, 6"
$t.Cell(9, 3).Range.Text  = "60÷8=7, 4"
$t.Cell(9, 4).Range.Text  = "86÷4=21, 2"
$t.Cell(9, 5).Range.Text  = "54÷3=18, 0"

$t.Cell(13, 1).Range.Text = "75÷7=10, 5"
$t.Cell(13, 2).Range.Text = "99÷5=19, 4"
$t.Cell(13, 3).Range.Text = "86÷8=10, 6"
$t.Cell(13, 4).Range.Text = "45÷5=9, 0"
$t.Cell(13, 5).Range.Text = "86÷7=12, 2"

$t.Cell(17, 1).Range.Text = "48÷7=6, 6"
$t.Cell(17, 2).Range.Text = "46÷4=11, 2"
$t.Cell(17, 3).Range.Text = "21÷3=7, 0"
$t.Cell(17, 4).Range.Text = "47÷3=15, 2"
$t.Cell(17, 5).Range.Text = "62÷6=10, 2"
